# Add data for 2022-11-10
# - Bump the "as of" date in the sheet name and the column B header text
#   from "November 01" to "November 02".
# - Add newly-recorded carjacking counts for a handful of
#   neighborhood/month cells (mostly historical "November" columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Through 2022-11-02"

# Update the column header text (row 1, column B) to match.
$ws.Range("B1").Value = "November 2022 (through November 02)"

# Englewood (row 2)
$ws.Range("M2").Value = 1
$ws.Range("AT2").Value = 1

# New City (row 4)
$ws.Range("M4").Value = 1
$ws.Range("BE4").Value = 1

# Garfield Park (row 5)
$ws.Range("B5").Value = 2
$ws.Range("M5").Value = 1
$ws.Range("X5").Value = 3

# Humboldt Park (row 7)
$ws.Range("B7").Value = 1
$ws.Range("X7").Value = 1

# Woodlawn (row 11)
$ws.Range("B11").Value = 1

# North Lawndale (row 14)
$ws.Range("M14").Value = 1
$ws.Range("X14").Value = 2

# Washington Heights (row 16)
$ws.Range("X16").Value = 1

# Ashburn (row 27)
$ws.Range("BP27").Value = 1

# Austin (row 28)
$ws.Range("B28").Value = 1
$ws.Range("AT28").Value = 2
$ws.Range("CA28").Value = 1

# West Ridge (row 32)
$ws.Range("M32").Value = 1

# West Loop (row 33)
$ws.Range("M33").Value = 1

# West Lawn (row 34)
$ws.Range("BP34").Value = 1

# Roseland (row 42)
$ws.Range("AI42").Value = 1
$ws.Range("BE42").Value = 1

# Brighton Park (row 56)
$ws.Range("M56").Value = 1

# Morgan Park (row 81)
$ws.Range("B81").Value = 1
$ws.Range("M81").Value = 2

# Near South Side (row 84)
$ws.Range("AT84").Value = 1

# Uptown (row 98)
$ws.Range("AI98").Value = 1
$ws.Range("AT98").Value = 2
